# Apply the "fixed example files" edit to the labelled.xlsx workbook:
#  - Append "/ concerned" to each of the 7 "...worried" response labels
#    (these live in column E "scare" and column F "overp", rows 2-8)
#  - Move the active selection to E2:E8
#  - Widen columns E:F to fit the new, longer label text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label text for the "worried / concerned" scale (rows 2-8 of columns E/F)
$labels = @(
    "Very, very worried/ concerned",
    "Very worried/ concerned",
    "Fairly worried/ concerned",
    "Somewhat worried/ concerned",
    "A little worried/ concerned",
    "Hardly worried/ concerned",
    "Not worried/ concerned"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $labels[$i]
    $ws.Cells.Item($row, 6).Value = $labels[$i]
}

# Widen columns E:F so the longer text fits (best-fit column width)
$ws.Columns("E:F").AutoFit()
$ws.Columns("E:F").ColumnWidth = 17.5

# Update the selection to match the edited range
$ws.Range("E2:E8").Select()
